$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume figures from the latest data pull.
# Column D ("Price") cells are forced to Text format before the value is set so
# that values like "1.00" or "62.658.21" stay text instead of being turned into
# numbers by Excel (matching the original inline-string cell type).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.658.21'
$ws.Range("E2").Value = '  -1.53%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.443.94'
$ws.Range("E3").Value = '  -1.51%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.26%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '569.82'
$ws.Range("E5").Value = '  -1.01%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.32'
$ws.Range("E6").Value = '  -3.95%  '
$ws.Range("E7").Value = '  +0.22%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.533'
$ws.Range("E8").Value = '  -1.58%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.439.53'
$ws.Range("E10").Value = '  -4.38%  '
$ws.Range("E11").Value = '  +1.33%  '
$ws.Range("E12").Value = '  -2.41%  '
$ws.Range("E13").Value = '  -2.67%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.06'
$ws.Range("E14").Value = '  -1.08%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000175'
$ws.Range("E15").Value = '  -5.65%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.886.96'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.539.37'
$ws.Range("E17").Value = '  -1.36%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.445.78'
$ws.Range("E18").Value = '  -1.76%  '
$ws.Range("E19").Value = '  -3.56%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.24'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '326.89'
$ws.Range("E21").Value = '  -0.82%  '
$ws.Range("E22").Value = '  -2.29%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.11'
$ws.Range("E23").Value = '  +10.59%  '
$ws.Range("E24").Value = '  +0.24%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '65.33'
$ws.Range("E25").Value = '  -3.63%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '626.23'
$ws.Range("E26").Value = '  -1.99%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.01'
$ws.Range("E27").Value = '  +1.50%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0₃0999'
$ws.Range("E28").Value = '  -5.88%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.565.96'
$ws.Range("E29").Value = '  -1.41%  '
$ws.Range("B30").Value = 'Binance-PegBSC-USD'
$ws.Range("C30").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  +0.57%  '
$ws.Range("B31").Value = 'Fetch.AI'
$ws.Range("C31").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.49'
$ws.Range("E31").Value = '  -3.07%  '
$ws.Range("E32").Value = '  -4.88%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.90'
$ws.Range("E33").Value = '  -1.67%  '
$ws.Range("E34").Value = '  -4.30%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.11'
$ws.Range("E35").Value = '  -2.15%  '
$ws.Range("E36").Value = '  -3.94%  '
$ws.Range("E37").Value = '  +0.26%  '
$ws.Range("E38").Value = '  -2.92%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.79'
$ws.Range("E39").Value = '  -1.02%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.32'
$ws.Range("E40").Value = '  -4.07%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '146.38'
$ws.Range("E41").Value = '  -0.91%  '
$ws.Range("E42").Value = '  -5.01%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.59'
$ws.Range("E43").Value = '  -2.32%  '
$ws.Range("E44").Value = '  +0.01%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '146.56'
$ws.Range("E45").Value = '  -3.76%  '
$ws.Range("E46").Value = '  -1.20%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '20.67'
$ws.Range("E47").Value = '  -2.54%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0529'
$ws.Range("E48").Value = '  -4.58%  '
$ws.Range("E49").Value = '  -2.55%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0232'
$ws.Range("E50").Value = '  -3.69%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0919'
$ws.Range("E51").Value = '  -1.23%  '
